{"js": "// Auto-generated replacement list: [before, after] pairs, in document order.\n// Each 'before' string is unique within the document, so a literal\n// whole-match search-and-replace is unambiguous for every entry.\nconst replacements = [\n  [\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"],\n  [\"18+22=40\", \"15+34=49\"],\n  [\"83-81=2\", \"15-12=3\"],\n  [\"2+42=44\", \"26+66=92\"],\n  [\"18+12=30\", \"83-0=83\"],\n  [\"87+2=89\", \"61+10=71\"],\n  [\"51-15=36\", \"70-39=31\"],\n  [\"28-23=5\", \"60+35=95\"],\n  [\"66-32=34\", \"52-47=5\"],\n  [\"77-2=75\", \"21+37=58\"],\n  [\"81-45=36\", \"87-87=0\"],\n  [\"6+53=59\", \"35+3=38\"],\n  [\"69+10=79\", \"67-21=46\"],\n  [\"97-42=55\", \"88-54=34\"],\n  [\"59+14=73\", \"25+2=27\"],\n  [\"40+55=95\", \"11+17=28\"],\n  [\"80-28=52\", \"94-40=54\"],\n  [\"78-35=43\", \"7-5=2\"],\n  [\"44-28=16\", \"41+31=72\"],\n  [\"19+26=45\", \"57+9=66\"],\n  [\"22-4=18\", \"69-18=51\"],\n  [\"12+24=36\", \"85-84=1\"],\n  [\"52-6=46\", \"94-29=65\"],\n  [\"69-57=12\", \"88-38=50\"],\n  [\"70-14=56\", \"41-11=30\"],\n  [\"11+65=76\", \"67+27=94\"],\n  [\"83-9=74\", \"88-60=28\"],\n  [\"67-37=30\", \"76+10=86\"],\n  [\"86-3=83\", \"32-26=6\"],\n  [\"64-31=33\", \"49+45=94\"],\n  [\"55-38=17\", \"75-31=44\"],\n  [\"27+21=48\", \"67+0=67\"],\n  [\"92-10=82\", \"67-32=35\"],\n  [\"64+3=67\", \"62+26=88\"],\n  [\"97-5=92\", \"90-41=49\"],\n  [\"49+4=53\", \"22+74=96\"],\n  [\"27+15=42\", \"77-51=26\"],\n  [\"5+59=64\", \"76-39=37\"],\n  [\"28+34=62\", \"6+20=26\"],\n  [\"25-9=16\", \"53+20=73\"],\n  [\"77-28=49\", \"45+41=86\"],\n  [\"51+16=67\", \"4+14=18\"],\n  [\"85-41=44\", \"83-55=28\"],\n  [\"98-96=2\", \"8+58=66\"],\n  [\"72-41=31\", \"27+45=72\"],\n  [\"90-72=18\", \"65-27=38\"],\n  [\"27+35=62\", \"8+56=64\"],\n  [\"56-46=10\", \"98-57=41\"],\n  [\"54-41=13\", \"69-42=27\"],\n  [\"98-76=22\", \"42-15=27\"],\n  [\"54-2=52\", \"22-3=19\"],\n  [\"61+22=83\", \"4+37=41\"],\n  [\"79-38=41\", \"43-0=43\"],\n  [\"65-46=19\", \"34+39=73\"],\n  [\"9+0=9\", \"33-12=21\"],\n  [\"79+18=97\", \"96-11=85\"],\n  [\"36+46=82\", \"24-20=4\"],\n  [\"68-64=4\", \"0+26=26\"],\n  [\"25+41=66\", \"16-8=8\"],\n  [\"58-36=22\", \"97-32=65\"],\n  [\"20+3=23\", \"9+81=90\"],\n  [\"94-14=80\", \"6+62=68\"],\n  [\"71-43=28\", \"20+15=35\"],\n  [\"61-2=59\", \"39-20=19\"],\n  [\"15+37=52\", \"33+40=73\"],\n  [\"26+25=51\", \"45+39=84\"],\n  [\"37-16=21\", \"52-26=26\"],\n  [\"53+34=87\", \"23+34=57\"],\n  [\"26+39=65\", \"71-16=55\"],\n  [\"17+24=41\", \"78+17=95\"],\n  [\"71-36=35\", \"94-91=3\"],\n  [\"33+19=52\", \"14+10=24\"],\n  [\"5+28=33\", \"23+32=55\"],\n  [\"39-28=11\", \"79-44=35\"],\n  [\"47+20=67\", \"24+37=61\"],\n  [\"32+43=75\", \"18+55=73\"],\n  [\"97-9=88\", \"8+0=8\"],\n  [\"37+15=52\", \"57+7=64\"],\n  [\"87-36=51\", \"3+77=80\"],\n  [\"51+19=70\", \"17+81=98\"],\n  [\"14+72=86\", \"61+20=81\"],\n  [\"22+63=85\", \"12+11=23\"],\n  [\"51-16=35\", \"3+72=75\"],\n  [\"81-79=2\", \"44-38=6\"],\n  [\"74+7=81\", \"2+68=70\"],\n  [\"55+36=91\", \"7+63=70\"],\n  [\"84-0=84\", \"58-39=19\"],\n  [\"36+4=40\", \"40-8=32\"],\n  [\"54+25=79\", \"97-75=22\"],\n  [\"6+31=37\", \"98-88=10\"],\n  [\"52-50=2\", \"24-13=11\"],\n  [\"63-1=62\", \"81-74=7\"],\n  [\"56-7=49\", \"26+32=58\"],\n  [\"56-35=21\", \"73+21=94\"],\n  [\"29+62=91\", \"78-24=54\"],\n  [\"53-13=40\", \"58-13=45\"],\n  [\"71-34=37\", \"2+6=8\"],\n  [\"16+25=41\", \"47-31=16\"],\n  [\"15-0=15\", \"28+20=48\"],\n  [\"23-20=3\", \"49-41=8\"],\n  [\"70-4=66\", \"38+27=65\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('text');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Auto-generated replacement list: (before, after) pairs, in document order.\n# Each 'before' string is unique within the document, so Find/Replace All\n# against the whole document content is unambiguous for every entry.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"2024-01-20 Saturday\", \"2024-01-21 Sunday\")\n  ,@(\"18+22=40\", \"15+34=49\")\n  ,@(\"83-81=2\", \"15-12=3\")\n  ,@(\"2+42=44\", \"26+66=92\")\n  ,@(\"18+12=30\", \"83-0=83\")\n  ,@(\"87+2=89\", \"61+10=71\")\n  ,@(\"51-15=36\", \"70-39=31\")\n  ,@(\"28-23=5\", \"60+35=95\")\n  ,@(\"66-32=34\", \"52-47=5\")\n  ,@(\"77-2=75\", \"21+37=58\")\n  ,@(\"81-45=36\", \"87-87=0\")\n  ,@(\"6+53=59\", \"35+3=38\")\n  ,@(\"69+10=79\", \"67-21=46\")\n  ,@(\"97-42=55\", \"88-54=34\")\n  ,@(\"59+14=73\", \"25+2=27\")\n  ,@(\"40+55=95\", \"11+17=28\")\n  ,@(\"80-28=52\", \"94-40=54\")\n  ,@(\"78-35=43\", \"7-5=2\")\n  ,@(\"44-28=16\", \"41+31=72\")\n  ,@(\"19+26=45\", \"57+9=66\")\n  ,@(\"22-4=18\", \"69-18=51\")\n  ,@(\"12+24=36\", \"85-84=1\")\n  ,@(\"52-6=46\", \"94-29=65\")\n  ,@(\"69-57=12\", \"88-38=50\")\n  ,@(\"70-14=56\", \"41-11=30\")\n  ,@(\"11+65=76\", \"67+27=94\")\n  ,@(\"83-9=74\", \"88-60=28\")\n  ,@(\"67-37=30\", \"76+10=86\")\n  ,@(\"86-3=83\", \"32-26=6\")\n  ,@(\"64-31=33\", \"49+45=94\")\n  ,@(\"55-38=17\", \"75-31=44\")\n  ,@(\"27+21=48\", \"67+0=67\")\n  ,@(\"92-10=82\", \"67-32=35\")\n  ,@(\"64+3=67\", \"62+26=88\")\n  ,@(\"97-5=92\", \"90-41=49\")\n  ,@(\"49+4=53\", \"22+74=96\")\n  ,@(\"27+15=42\", \"77-51=26\")\n  ,@(\"5+59=64\", \"76-39=37\")\n  ,@(\"28+34=62\", \"6+20=26\")\n  ,@(\"25-9=16\", \"53+20=73\")\n  ,@(\"77-28=49\", \"45+41=86\")\n  ,@(\"51+16=67\", \"4+14=18\")\n  ,@(\"85-41=44\", \"83-55=28\")\n  ,@(\"98-96=2\", \"8+58=66\")\n  ,@(\"72-41=31\", \"27+45=72\")\n  ,@(\"90-72=18\", \"65-27=38\")\n  ,@(\"27+35=62\", \"8+56=64\")\n  ,@(\"56-46=10\", \"98-57=41\")\n  ,@(\"54-41=13\", \"69-42=27\")\n  ,@(\"98-76=22\", \"42-15=27\")\n  ,@(\"54-2=52\", \"22-3=19\")\n  ,@(\"61+22=83\", \"4+37=41\")\n  ,@(\"79-38=41\", \"43-0=43\")\n  ,@(\"65-46=19\", \"34+39=73\")\n  ,@(\"9+0=9\", \"33-12=21\")\n  ,@(\"79+18=97\", \"96-11=85\")\n  ,@(\"36+46=82\", \"24-20=4\")\n  ,@(\"68-64=4\", \"0+26=26\")\n  ,@(\"25+41=66\", \"16-8=8\")\n  ,@(\"58-36=22\", \"97-32=65\")\n  ,@(\"20+3=23\", \"9+81=90\")\n  ,@(\"94-14=80\", \"6+62=68\")\n  ,@(\"71-43=28\", \"20+15=35\")\n  ,@(\"61-2=59\", \"39-20=19\")\n  ,@(\"15+37=52\", \"33+40=73\")\n  ,@(\"26+25=51\", \"45+39=84\")\n  ,@(\"37-16=21\", \"52-26=26\")\n  ,@(\"53+34=87\", \"23+34=57\")\n  ,@(\"26+39=65\", \"71-16=55\")\n  ,@(\"17+24=41\", \"78+17=95\")\n  ,@(\"71-36=35\", \"94-91=3\")\n  ,@(\"33+19=52\", \"14+10=24\")\n  ,@(\"5+28=33\", \"23+32=55\")\n  ,@(\"39-28=11\", \"79-44=35\")\n  ,@(\"47+20=67\", \"24+37=61\")\n  ,@(\"32+43=75\", \"18+55=73\")\n  ,@(\"97-9=88\", \"8+0=8\")\n  ,@(\"37+15=52\", \"57+7=64\")\n  ,@(\"87-36=51\", \"3+77=80\")\n  ,@(\"51+19=70\", \"17+81=98\")\n  ,@(\"14+72=86\", \"61+20=81\")\n  ,@(\"22+63=85\", \"12+11=23\")\n  ,@(\"51-16=35\", \"3+72=75\")\n  ,@(\"81-79=2\", \"44-38=6\")\n  ,@(\"74+7=81\", \"2+68=70\")\n  ,@(\"55+36=91\", \"7+63=70\")\n  ,@(\"84-0=84\", \"58-39=19\")\n  ,@(\"36+4=40\", \"40-8=32\")\n  ,@(\"54+25=79\", \"97-75=22\")\n  ,@(\"6+31=37\", \"98-88=10\")\n  ,@(\"52-50=2\", \"24-13=11\")\n  ,@(\"63-1=62\", \"81-74=7\")\n  ,@(\"56-7=49\", \"26+32=58\")\n  ,@(\"56-35=21\", \"73+21=94\")\n  ,@(\"29+62=91\", \"78-24=54\")\n  ,@(\"53-13=40\", \"58-13=45\")\n  ,@(\"71-34=37\", \"2+6=8\")\n  ,@(\"16+25=41\", \"47-31=16\")\n  ,@(\"15-0=15\", \"28+20=48\")\n  ,@(\"23-20=3\", \"49-41=8\")\n  ,@(\"70-4=66\", \"38+27=65\")\n)\n\nforeach ($pair in $replacements) {\n  $before = $pair[0]\n  $after = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $found = $range.Find.Execute($before, $true, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n  if (-not $found) {\n    throw \"No match found for: $before\"\n  }\n}"}
